$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns E:F (SLAND_corrected, SLAND_dor)
$ws.Range("E1:F1").EntireColumn.Delete() | Out-Null

# Remove columns J:M (F_ab_cor, F_abc_cor, F_ab_dor, F_abc_dor) after the shift
$ws.Range("J1:M1").EntireColumn.Delete() | Out-Null

# Update header row and data values for the new layout (A:K)
$ws.Range("A1").Value = 'continent'
$ws.Range("B1").Value = 'iso'
$ws.Range("C1").Value = 'CSCC'
$ws.Range("D1").Value = 'SLAND'
$ws.Range("E1").Value = 'ELUC'
$ws.Range("F1").Value = 'EFOS'
$ws.Range("G1").Value = 'F_ab'
$ws.Range("H1").Value = 'F_ac'
$ws.Range("I1").Value = 'F_abc'
$ws.Range("J1").Value = 'economic_group'
$ws.Range("K1").Value = 'climate_zone'
$ws.Range("A2").Value = 'Africa'
$ws.Range("B2").Value = 'AGOBDIBENBFABWACAFCIVCMRCODCOGCOMCPVDJIDZAEGYERIETHGABGHAGINGMBGNBGNQKENLBRLBYLSOMARMDGMLIMOZMRTMUSMWINAMNERNGARWASDNSENSLESOMSOMSSDSTPSWZTCDTGOTUNTZAUGAZAFZMBZWE'
$ws.Range("C2").Value = 176.3343492289424
$ws.Range("D2").Value = 0.7245317697525024
$ws.Range("E2").Value = -0.37219005425
$ws.Range("F2").Value = -0.357894456830619
$ws.Range("G2").Value = 0.3523417224354953
$ws.Range("H2").Value = -0.730084511080619
$ws.Range("I2").Value = -0.00555273439512376
$ws.Range("J2").Value = 'LDCLDCLDCLDCOtherLDCOtherOtherLDCOtherLDCOtherLDCOtherBRICSLDCBRICSOtherOtherLDCLDCLDCOtherOtherLDCOtherLDCOtherLDCLDCLDCOtherOtherLDCOtherLDCOtherLDCOtherLDCLDCOtherOtherLDCLDCLDCLDCLDCOtherLDCLDCBRICSLDCLDC'
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = '111121111112222211111111122212121122112112221121211212'
$ws.Range("K2").ClearFormats()
$ws.Range("A3").Value = 'Asia'
$ws.Range("B3").Value = 'AFGAREARMAZEBGDBRNBTNCHNCYPGEOIDNINDIRNIRQISRJORJPNKAZKGZKHMKORKWTLAOLBNLKAMMRMNGMYSNPLOMNPAKPHLPRKPSEQATRUSSAUSYRTHATJKTKMTLSTURUZBVNMYEM'
$ws.Range("C3").Value = 409.9114079969629
$ws.Range("D3").Value = 1.23817253112793
$ws.Range("E3").Value = -0.36132231
$ws.Range("F3").Value = -5.809045108948109
$ws.Range("G3").Value = 0.8768502678982399
$ws.Range("H3").Value = -6.170367418948108
$ws.Range("I3").Value = -4.932194841049868
$ws.Range("J3").Value = 'LDCBRICSOtherOtherOtherOtherOtherBRICSOtherOtherOtherBRICSBRICSOtherOECDOtherOECDOtherOtherLDCOECDOtherLDCOtherOtherLDCOtherOtherOtherOtherOtherOtherOtherOtherOtherBRICSOtherOtherOtherOtherOtherLDCOECDOtherOtherOther'
$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = '4244114434112222344132131141422142252214214412'
$ws.Range("K3").ClearFormats()
$ws.Range("A4").Value = 'Europe'
$ws.Range("B4").Value = 'ALBAUTBELBGRBIHBLRCHECZEDEUDNKESPESTFINFRAGBRGRCHRVHUNIRLISLITALTULUXLVAMDAMKDMNENLDNORPOLPRTROUSRBSVKSVNSWEUKR'
$ws.Range("C4").Value = -19.44382776680462
$ws.Range("D4").Value = 0.1197896003723145
$ws.Range("E4").Value = 0.02725090325
$ws.Range("F4").Value = -1.035789632326314
$ws.Range("G4").Value = 0.1470405011121091
$ws.Range("H4").Value = -1.008538729076314
$ws.Range("I4").Value = -0.8887491312142051
$ws.Range("J4").Value = 'OtherOECDOECDOtherOtherOtherOECDOECDOECDOECDOECDOECDOECDOECDOECDOECDOtherOECDOECDOECDOECDOECDOECDOECDOtherOtherOtherOECDOECDOECDOECDOtherOtherOECDOECDOECDOther'
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = '3333343333344333333534344333543333344'
$ws.Range("K4").ClearFormats()
$ws.Range("A5").Value = 'Latin America and the Caribbean'
$ws.Range("B5").Value = 'ARGBHSBLZBOLBRACHLCOLCRICUBDOMECUGTMGUYHNDHTIJAMMEXNICPANPERPRYSLVSURTTOURYVCTVEN'
$ws.Range("C5").Value = 69.63392779362381
$ws.Range("D5").Value = 0.6940719485282898
$ws.Range("E5").Value = -0.4263971400000001
$ws.Range("F5").Value = -0.4686550021483723
$ws.Range("G5").Value = 0.2676748212739485
$ws.Range("H5").Value = -0.8950521421483724
$ws.Range("I5").Value = -0.2009801808744238
$ws.Range("J5").Value = 'OtherOtherOtherOtherBRICSOECDOECDOECDOtherOtherOtherOtherOtherOtherLDCOtherOECDOtherOtherOtherOtherOtherOtherOtherOtherOtherOther'
$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = '111112111111111111111111111'
$ws.Range("K5").ClearFormats()
$ws.Range("A6").Value = 'North America'
$ws.Range("B6").Value = 'CANUSA'
$ws.Range("C6").Value = 21.59084712077145
$ws.Range("D6").Value = 0.5105834603309631
$ws.Range("E6").Value = -0.02196698900000001
$ws.Range("F6").Value = -1.580879659226444
$ws.Range("G6").Value = 0.4886164564298019
$ws.Range("H6").Value = -1.602846648226444
$ws.Range("I6").Value = -1.092263202796642
$ws.Range("J6").Value = 'OECDOECD'
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = '53'
$ws.Range("K6").ClearFormats()
$ws.Range("A7").Value = 'Oceania'
$ws.Range("B7").Value = 'AUSCOKFJIKIRNZLPNGSLBVUTWSM'
$ws.Range("C7").Value = 4.413051608622875
$ws.Range("D7").Value = 0.08273434638977051
$ws.Range("E7").Value = -0.01988379875
$ws.Range("F7").Value = -0.1217305771283392
$ws.Range("G7").Value = 0.06285054773168632
$ws.Range("H7").Value = -0.1416143758783392
$ws.Range("I7").Value = -0.05888002939665294
$ws.Range("J7").Value = 'OECDOtherOtherOtherOECDOtherOtherOtherOther'
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = '211131111'
$ws.Range("K7").ClearFormats()
